$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B, shifting the existing columns
# B..O to C..P. This mirrors the "Internal" column that was added both
# to the SAMPLE_TYPE table (row 2/3) and to the property-assignment
# table (row 4-9).
$ws.Columns("B:B").Insert()

# --- SAMPLE_TYPE table header (row 2) ---
$ws.Range("B2").Value = "Internal"

# --- SAMPLE_TYPE table data (row 3) ---
# Leading apostrophe forces the "TRUE"/"FALSE" text (matching the rest
# of the sheet) instead of Excel auto-coercing it to a boolean.
$ws.Range("B3").Value = "'FALSE"

# --- Property assignment table header (row 4) ---
$ws.Range("B4").Value = "Internal"

# Append the new "Internal Assignment" column at the end of the
# property-assignment table.
$ws.Range("P4").Value = "Internal Assignment"

# --- Property assignment table data (rows 5-9) ---
$ws.Range("B5").Value = "'TRUE"
$ws.Range("B6").Value = "'FALSE"
$ws.Range("B7").Value = "'FALSE"
$ws.Range("B8").Value = "'FALSE"
$ws.Range("B9").Value = "'FALSE"

$ws.Range("P5").Value = "'FALSE"
$ws.Range("P6").Value = "'FALSE"
$ws.Range("P7").Value = "'FALSE"
$ws.Range("P8").Value = "'FALSE"
$ws.Range("P9").Value = "'FALSE"

# Fix the NAME property code (was "$NAME").
$ws.Range("A5").Value = "NAME"

# Restore the style used for the new Internal Assignment rows (matches
# the boolean-like TRUE/FALSE display style used elsewhere in the
# property table).
$ws.Range("P5:P9").Style = $ws.Range("L5").Style

# Match the selection left behind by the column insert in the source
# workbook.
$ws.Range("B2:B6").Select()
